# Add a "Quantity" column to the Product sheet, populate it with 10 units
# for every existing product row, and leave the workbook with the Product
# sheet active (selection on F4), matching the "add quantity to database"
# commit.

$wb = $excel.ActiveWorkbook

$wsCategory = $wb.Worksheets.Item("Category")
$wsProduct  = $wb.Worksheets.Item("Product")

# --- Product sheet: new "Quantity" column (E) ---------------------------
# Match the formatting (thin border style) already used by the rest of
# the header/data rows on this sheet before filling in values.
$wsProduct.Range("D1").Copy()
$wsProduct.Range("E1:E10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsProduct.Range("E1").Value = "Quantity"

$quantities = @(10, 10, 10, 10, 10, 10, 10, 10, 10)
for ($i = 0; $i -lt $quantities.Length; $i++) {
    $row = $i + 2
    $wsProduct.Cells.Item($row, 5).Value = $quantities[$i]
}

# --- Selections / active sheet ------------------------------------------
# Category sheet is no longer the active tab; its lingering selection
# moves to B1.
$wsCategory.Range("B1").Select()

# Product becomes the active sheet, with the selection parked on F4.
$wsProduct.Activate()
$wsProduct.Range("F4").Select()
